# Add 10 homes for non-solar homes
# Copy the "Monthly 10 homes" pivot block (B1:L16) and paste it next to the
# existing pivot table on "Monthly 10 homes 2" (landing at N1:X16).

$wb = $excel.ActiveWorkbook

$wsSource = $wb.Worksheets.Item("Monthly 10 homes")
$wsDest   = $wb.Worksheets.Item("Monthly 10 homes 2")

# Copy source range B1:L16 from the "Monthly 10 homes" sheet
$wsSource.Activate()
$srcRange = $wsSource.Range("B1:L16")
$srcRange.Copy()

# Paste onto the "Monthly 10 homes 2" sheet starting at N1
$wsDest.Activate()
$destCell = $wsDest.Range("N1")
$destCell.PasteSpecial(-4104)

$excel.CutCopyMode = 0

# Update selection/view state to mirror post-copy/paste state
$wsSource.Range("A1:L16").Select()
$wsSource.Application.ActiveWindow.ScrollRow = 1

$wsDest.Range("Q26").Select()
